$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap columns F:V between row pairs (home/away fixtures were reordered) ---
# rows 14 <-> 15
$ws.Range("F14").Value = "Benfica"
$ws.Range("F15").Value = "Portimonense"
$ws.Range("G14").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H14").Value = "Estrela"
$ws.Range("H15").Value = "Boavista"
$ws.Range("I14").Value = 0
$ws.Range("I15").Value = 4
$ws.Range("J14").Value = 1.12
$ws.Range("J15").Value = 2.45
$ws.Range("K14").Value = "15/08/2023 11:42"
$ws.Range("K15").Value = "15/08/2023 11:42"
$ws.Range("L14").Value = 1.15
$ws.Range("L15").Value = 3.35
$ws.Range("M14").Value = "19/08/2023 20:57"
$ws.Range("M15").Value = "19/08/2023 21:28"
$ws.Range("N14").Value = 10.69
$ws.Range("N15").Value = 3.31
$ws.Range("O14").Value = "15/08/2023 11:42"
$ws.Range("O15").Value = "15/08/2023 11:42"
$ws.Range("P14").Value = 8.859999999999999
$ws.Range("P15").Value = 3.24
$ws.Range("Q14").Value = "19/08/2023 21:25"
$ws.Range("Q15").Value = "19/08/2023 21:22"
$ws.Range("R14").Value = 22.27
$ws.Range("R15").Value = 3.16
$ws.Range("S14").Value = "15/08/2023 11:42"
$ws.Range("S15").Value = "15/08/2023 11:42"
$ws.Range("T14").Value = 19.04
$ws.Range("T15").Value = 2.38
$ws.Range("U14").Value = "19/08/2023 21:25"
$ws.Range("U15").Value = "19/08/2023 21:28"
$ws.Range("V14").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/benfica-estrela-da-amadora/zFdna79Q/"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/portimonense-boavista/C0wwkNoJ/"

# rows 66 <-> 67
$ws.Range("F66").Value = "SC Farense"
$ws.Range("F67").Value = "Chaves"
$ws.Range("G66").Value = 0
$ws.Range("G67").Value = 4
$ws.Range("H66").Value = "Vizela"
$ws.Range("H67").Value = "Gil Vicente"
$ws.Range("I66").Value = 0
$ws.Range("I67").Value = 2
$ws.Range("J66").Value = 2.29
$ws.Range("J67").Value = 2.62
$ws.Range("K66").Value = "02/10/2023 07:12"
$ws.Range("K67").Value = "02/10/2023 20:42"
$ws.Range("L66").Value = 2.32
$ws.Range("L67").Value = 2.74
$ws.Range("M66").Value = "07/10/2023 16:02"
$ws.Range("M67").Value = "07/10/2023 16:29"
$ws.Range("N66").Value = 3.4
$ws.Range("N67").Value = 3.45
$ws.Range("O66").Value = "02/10/2023 07:12"
$ws.Range("O67").Value = "02/10/2023 20:42"
$ws.Range("P66").Value = 3.47
$ws.Range("P67").Value = 3.65
$ws.Range("Q66").Value = "07/10/2023 15:49"
$ws.Range("Q67").Value = "07/10/2023 16:25"
$ws.Range("R66").Value = 3.31
$ws.Range("R67").Value = 2.73
$ws.Range("S66").Value = "02/10/2023 07:12"
$ws.Range("S67").Value = "02/10/2023 20:42"
$ws.Range("T66").Value = 3.24
$ws.Range("T67").Value = 2.58
$ws.Range("U66").Value = "07/10/2023 16:03"
$ws.Range("U67").Value = "07/10/2023 16:29"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/sc-farense-vizela/OY1Asc0E/"
$ws.Range("V67").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/chaves-gil-vicente/K4BKKZh1/"

# rows 70 <-> 71
$ws.Range("F70").Value = "Casa Pia"
$ws.Range("F71").Value = "Famalicao"
$ws.Range("G70").Value = 0
$ws.Range("G71").Value = 1
$ws.Range("H70").Value = "Estrela"
$ws.Range("H71").Value = "Vitoria Guimaraes"
$ws.Range("I70").Value = 1
$ws.Range("I71").Value = 3
$ws.Range("J70").Value = 1.89
$ws.Range("J71").Value = 2.59
$ws.Range("K70").Value = "02/10/2023 20:42"
$ws.Range("K71").Value = "02/10/2023 07:12"
$ws.Range("L70").Value = 2.16
$ws.Range("L71").Value = 2.94
$ws.Range("M70").Value = "08/10/2023 16:29"
$ws.Range("M71").Value = "08/10/2023 16:27"
$ws.Range("N70").Value = 3.53
$ws.Range("N71").Value = 3.26
$ws.Range("O70").Value = "02/10/2023 20:42"
$ws.Range("O71").Value = "02/10/2023 07:12"
$ws.Range("P70").Value = 3.47
$ws.Range("P71").Value = 3.07
$ws.Range("Q70").Value = "08/10/2023 16:29"
$ws.Range("Q71").Value = "08/10/2023 16:27"
$ws.Range("R70").Value = 4.36
$ws.Range("R71").Value = 2.89
$ws.Range("S70").Value = "02/10/2023 20:42"
$ws.Range("S71").Value = "02/10/2023 07:12"
$ws.Range("T70").Value = 3.61
$ws.Range("T71").Value = 2.77
$ws.Range("U70").Value = "08/10/2023 16:28"
$ws.Range("U71").Value = "08/10/2023 16:27"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/casa-pia-estrela-da-amadora/Cbb6rwo8/"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/famalicao-vitoria-guimaraes/beAOJg87/"

# rows 76 <-> 77
$ws.Range("F76").Value = "Vitoria Guimaraes"
$ws.Range("F77").Value = "Benfica"
$ws.Range("G76").Value = 5
$ws.Range("G77").Value = 1
$ws.Range("H76").Value = "Chaves"
$ws.Range("H77").Value = "Casa Pia"
$ws.Range("I76").Value = 0
$ws.Range("I77").Value = 1
$ws.Range("J76").Value = 1.62
$ws.Range("J77").Value = 1.22
$ws.Range("K76").Value = "11/10/2023 14:42"
$ws.Range("K77").Value = "11/10/2023 14:42"
$ws.Range("L76").Value = 1.69
$ws.Range("L77").Value = 1.22
$ws.Range("M76").Value = "28/10/2023 18:58"
$ws.Range("M77").Value = "28/10/2023 18:55"
$ws.Range("N76").Value = 4.26
$ws.Range("N77").Value = 7.32
$ws.Range("O76").Value = "11/10/2023 14:42"
$ws.Range("O77").Value = "11/10/2023 14:42"
$ws.Range("P76").Value = 4
$ws.Range("P77").Value = 6.95
$ws.Range("Q76").Value = "28/10/2023 18:58"
$ws.Range("Q77").Value = "28/10/2023 18:58"
$ws.Range("R76").Value = 5.63
$ws.Range("R77").Value = 13.22
$ws.Range("S76").Value = "11/10/2023 14:42"
$ws.Range("S77").Value = "11/10/2023 14:42"
$ws.Range("T76").Value = 5.29
$ws.Range("T77").Value = 14.17
$ws.Range("U76").Value = "28/10/2023 18:58"
$ws.Range("U77").Value = "28/10/2023 18:58"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/vitoria-guimaraes-chaves/8vH9wlat/"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/benfica-casa-pia/GWtkzFhl/"

# --- Append new rows 93-95 (copy formatting from row 92, then fill values) ---
$ws.Range("A92:V92").Copy()
$ws.Range("A93:V95").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 93
$ws.Range("A93").Value = 92
$ws.Range("B93").Value = "portugal"
$ws.Range("C93").Value = "liga-portugal"
$ws.Range("D93").Value = "2023-2024"
$ws.Range("E93").Value = 45241.6875
$ws.Range("F93").Value = "Portimonense"
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = "Chaves"
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = 2.12
$ws.Range("K93").Value = "05/11/2023 15:44"
$ws.Range("L93").Value = 2.23
$ws.Range("M93").Value = "11/11/2023 16:27"
$ws.Range("N93").Value = 3.68
$ws.Range("O93").Value = "05/11/2023 15:44"
$ws.Range("P93").Value = 3.6
$ws.Range("Q93").Value = "11/11/2023 16:28"
$ws.Range("R93").Value = 3.46
$ws.Range("S93").Value = "05/11/2023 15:44"
$ws.Range("T93").Value = 3.31
$ws.Range("U93").Value = "11/11/2023 16:28"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/portimonense-chaves/UyCnu1zo/"

# row 94
$ws.Range("A94").Value = 93
$ws.Range("B94").Value = "portugal"
$ws.Range("C94").Value = "liga-portugal"
$ws.Range("D94").Value = "2023-2024"
$ws.Range("E94").Value = 45241.79166666666
$ws.Range("F94").Value = "Estrela"
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = "Moreirense"
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 2.58
$ws.Range("K94").Value = "06/11/2023 03:12"
$ws.Range("L94").Value = 2.77
$ws.Range("M94").Value = "11/11/2023 18:58"
$ws.Range("N94").Value = 3.36
$ws.Range("O94").Value = "06/11/2023 03:12"
$ws.Range("P94").Value = 3.51
$ws.Range("Q94").Value = "11/11/2023 18:54"
$ws.Range("R94").Value = 2.83
$ws.Range("S94").Value = "06/11/2023 03:12"
$ws.Range("T94").Value = 2.63
$ws.Range("U94").Value = "11/11/2023 18:58"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/estrela-da-amadora-moreirense/SSPKqsST/"

# row 95
$ws.Range("A95").Value = 94
$ws.Range("B95").Value = "portugal"
$ws.Range("C95").Value = "liga-portugal"
$ws.Range("D95").Value = "2023-2024"
$ws.Range("E95").Value = 45241.79166666666
$ws.Range("F95").Value = "Vizela"
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = "Famalicao"
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 2.35
$ws.Range("K95").Value = "05/11/2023 16:42"
$ws.Range("L95").Value = 2.42
$ws.Range("M95").Value = "11/11/2023 18:59"
$ws.Range("N95").Value = 3.36
$ws.Range("O95").Value = "05/11/2023 16:42"
$ws.Range("P95").Value = 3.28
$ws.Range("Q95").Value = "11/11/2023 18:31"
$ws.Range("R95").Value = 3.23
$ws.Range("S95").Value = "05/11/2023 16:42"
$ws.Range("T95").Value = 3.24
$ws.Range("U95").Value = "11/11/2023 18:59"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/portugal/liga-portugal/vizela-famalicao/WtjV3pDo/"

